$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.720.30"
$ws.Range("E2").Value = "  -0.30%  "

$ws.Range("D3").Value = "'1.633.21"
$ws.Range("E3").Value = "  -0.16%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'215.24"
$ws.Range("E5").Value = "  +0.10%  "

$ws.Range("E6").Value = "  -0.62%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").Value = "'0.0635"
$ws.Range("E9").Value = "  -1.02%  "

$ws.Range("D10").Value = "'19.59"
$ws.Range("E10").Value = "  -3.34%  "

$ws.Range("E11").Value = "  +0.92%  "

$ws.Range("E12").Value = "  -0.28%  "

$ws.Range("D13").Value = "'1.859.48"
$ws.Range("E13").Value = "  -0.12%  "

$ws.Range("D14").Value = "'1.634.15"
$ws.Range("E14").Value = "  -0.09%  "

$ws.Range("D15").Value = "'0.556"
$ws.Range("E15").Value = "  -0.62%  "

$ws.Range("E16").Value = "  -0.14%  "

$ws.Range("D17").Value = "'62.68"
$ws.Range("E17").Value = "  -0.91%  "

$ws.Range("D18").Value = "'25.754.32"
$ws.Range("E18").Value = "  -0.21%  "

$ws.Range("E19").Value = "  -0.02%  "

$ws.Range("D20").Value = "'4.44"
$ws.Range("E20").Value = "  +1.52%  "

$ws.Range("D21").Value = "'193.74"
$ws.Range("E21").Value = "  +0.81%  "

$ws.Range("D22").Value = "'9.93"
$ws.Range("E22").Value = "  +0.30%  "

$ws.Range("D23").Value = "'6.28"
$ws.Range("E23").Value = "  +2.29%  "

$ws.Range("D25").Value = "'1.83"
$ws.Range("E25").Value = "  +3.42%  "

$ws.Range("D26").Value = "'141.84"
$ws.Range("E26").Value = "  +2.10%  "

$ws.Range("E27").Value = "  -1.10%  "

$ws.Range("D28").Value = "'6.88"
$ws.Range("E28").Value = "  +0.99%  "

$ws.Range("D29").Value = "'15.50"
$ws.Range("E29").Value = "  -0.11%  "

$ws.Range("E30").Value = "  -0.06%  "

$ws.Range("D31").Value = "'0.0491"
$ws.Range("E31").Value = "  -0.47%  "

$ws.Range("D32").Value = "'3.32"
$ws.Range("E32").Value = "  +1.03%  "

$ws.Range("D33").Value = "'3.24"
$ws.Range("E33").Value = "  -0.36%  "

$ws.Range("E34").Value = "  +0.39%  "

$ws.Range("E35").Value = "  +0.18%  "

$ws.Range("D36").Value = "'0.900"
$ws.Range("E36").Value = "  -0.34%  "

$ws.Range("D37").Value = "'1.123.28"
$ws.Range("E37").Value = "  -0.70%  "

$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "'2.53"
$ws.Range("E38").Value = "  -1.59%  "

$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "'0.547"
$ws.Range("E39").Value = "  -1.59%  "

$ws.Range("E40").Value = "  -1.02%  "

$ws.Range("E41").Value = "  +0.70%  "

$ws.Range("D42").Value = "'5.57"
$ws.Range("E42").Value = "  +2.19%  "

$ws.Range("D43").Value = "'99.71"
$ws.Range("E43").Value = "  +0.90%  "

$ws.Range("D44").Value = "'0.803"
$ws.Range("E44").Value = "  +0.66%  "

$ws.Range("D45").Value = "'1.768.16"
$ws.Range("E45").Value = "  -0.24%  "

$ws.Range("E46").Value = "  +2.39%  "

$ws.Range("D47").Value = "'55.01"
$ws.Range("E47").Value = "  -1.05%  "

$ws.Range("E49").Value = "  -0.11%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'7.56"
$ws.Range("E50").Value = "  -2.97%  "

$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").Value = "'2.33"
$ws.Range("E51").Value = "  +2.68%  "
